# Updated cryptos list values (Price / Volume(1h)) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether Excel would
# auto-convert the text to a number (single-decimal "price" cells like
# "30.74") and so needs to be forced back to plain text, matching the
# source data's text-typed cells (e.g. multi-dot thousands values like
# "33.938.51" already survive as text on their own).
$updates = @(
    @{ Cell = 'D2'; Value = '33.938.51'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +9.17%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.783.14'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +5.76%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.21%  '; ForceText = $false }
    @{ Cell = 'E5'; Value = '  +2.05%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  +4.08%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.30%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '30.74'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +4.58%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '46.42'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +3.02%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.279'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +4.46%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  +2.99%  '; ForceText = $false }
    @{ Cell = 'E12'; Value = '  +1.94%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '2.041.50'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  +5.83%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '1.782.81'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  +5.99%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.628'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +3.02%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '33.915.70'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +8.86%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '9.99'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -3.68%  '; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +1.05%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '68.59'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +2.69%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '251.52'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +0.52%  '; ForceText = $false }
    @{ Cell = 'E21'; Value = '  +2.60%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  +0.23%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '10.30'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +2.51%  '; ForceText = $false }
    @{ Cell = 'E24'; Value = '  -1.99%  '; ForceText = $false }
    @{ Cell = 'E25'; Value = '  -0.25%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '158.43'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '16.47'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +3.27%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +1.60%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '6.95'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +3.39%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  +0.29%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  +7.59%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '0.0513'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +2.81%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  +3.94%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '3.55'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +6.13%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.492.29'; ForceText = $false }
    @{ Cell = 'E35'; Value = '  -1.59%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '1.79'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +2.60%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  +3.02%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  +3.22%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  +2.61%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '83.15'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'E41'; Value = '  +3.04%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '2.72'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +0.31%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.887'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +5.23%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  +2.05%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.0509'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +1.11%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '1.07'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +3.15%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '1.938.77'; ForceText = $false }
    @{ Cell = 'E47'; Value = '  +6.48%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '5.78'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +3.59%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  +0.35%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '11.91'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +15.21%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '50.89'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -2.76%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $r.NumberFormat = "@"
        $r.Value = $u.Value
        $r.ClearFormats()
    } else {
        $r.Value = $u.Value
    }
}
